$wb = $excel.ActiveWorkbook

# Target row-2 values per sheet, indexed by column letter (after insertion of gb + btes columns)
$values = @{
    1 = @{ "A"=0; "B"=0; "C"=34850.02913702199; "D"=0; "E"=695202.8899878451; "F"=1000.138367865942; "G"=0; "H"=2534.277928792126; "I"=0; "J"=0; "K"=0; "L"=0; "M"=0; "N"=218031.1156391199; "O"=1996.112248849872 }
    2 = @{ "A"=0; "B"=0; "C"=149504.1992614464; "D"=0; "E"=695202.8899878451; "F"=1000.138367865942; "G"=0; "H"=7448.027758327713; "I"=0; "J"=0; "K"=0; "L"=0; "M"=0; "N"=228053.6181928575; "O"=5709.431641187627 }
    3 = @{ "A"=0; "B"=0; "C"=253633.6218973605; "D"=0; "E"=695202.8899878451; "F"=15729.44046117129; "G"=0; "H"=11595.71584387759; "I"=0; "J"=0; "K"=0; "L"=0; "M"=0; "N"=233561.946346158; "O"=9655.767294985335 }
    4 = @{ "A"=0; "B"=0; "C"=253633.6218973605; "D"=0; "E"=695202.8899878451; "F"=15729.44046117129; "G"=0; "H"=11595.71584387759; "I"=0; "J"=0; "K"=0; "L"=0; "M"=0; "N"=233561.9463461617; "O"=9655.767294985335 }
    5 = @{ "A"=0; "B"=0; "C"=253633.6218973605; "D"=0; "E"=695202.8899878451; "F"=15729.44046117129; "G"=0; "H"=11595.71584387759; "I"=0; "J"=0; "K"=0; "L"=0; "M"=0; "N"=233561.9463461617; "O"=9655.767294985335 }
    6 = @{ "A"=0; "B"=0; "C"=253633.6218973605; "D"=0; "E"=695202.8899878451; "F"=15729.44046117129; "G"=0; "H"=11595.71584387759; "I"=0; "J"=0; "K"=0; "L"=0; "M"=0; "N"=233561.9463461617; "O"=9655.767294985335 }
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Insert new column B for "gb" (shifts hp..ites from B..M to C..N)
    $ws.Columns.Item(2).Insert()
    $ws.Range("B1").Value = "gb"

    # Insert new column N for "btes" (shifts ites from N to O)
    $ws.Columns.Item(14).Insert()
    $ws.Range("N1").Value = "btes"

    $rowVals = $values[$i]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col" + "2").Value = $rowVals[$col]
    }
}
